# #5: insurance, claim, debt, investment done
#
#  - "具有相當價值之財產" sheet: the antiques/jewelry rows had been tagged
#    with the wrong property_category ("otherbonds"); relabel as "antique".
#  - "保險" (insurance) sheet: re-export with the full normalized schema
#    (company, name, owner, property_category, category, date,
#     legislator_name, legislator_id, source_file, index) and fix the
#     per-row owner/category values.
#  - "債務" (debt) sheet: same normalization (species, debtor, owner,
#     total, register_date, register_reason, property_category, category,
#     date, legislator_name, legislator_id, source_file, index); the old
#     header row had stray data values instead of real column names.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 0. Fix the mislabeled shared string used by the "具有相當價值之財產"
#    sheet (antiques/jewelry rows were tagged "otherbonds").
# ------------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("具有相當價值之財產")
$wsAssets.Range("F2").Value = "antique"
$wsAssets.Range("F3").Value = "antique"

# ------------------------------------------------------------------
# 1. "保險" (insurance) sheet
# ------------------------------------------------------------------
$wsIns = $wb.Worksheets.Item("保險")

# Header row
$wsIns.Range("B1").Value = "company"
$wsIns.Range("C1").Value = "name"
$wsIns.Range("D1").Value = "owner"
$wsIns.Range("E1").Value = "property_category"
$wsIns.Range("F1").Value = "category"
$wsIns.Range("G1").Value = "date"
$wsIns.Range("H1").Value = "legislator_name"
$wsIns.Range("I1").Value = "legislator_id"
$wsIns.Range("J1").Value = "source_file"
$wsIns.Range("K1").Value = "index"

# Row 2 (index 105)
$wsIns.Range("B2").Value = "南山人壽"
$wsIns.Range("C2").Value = "南山康福20年期終身壽險"
$wsIns.Range("D2").Value = "丁守中"
$wsIns.Range("E2").Value = "insurance"
$wsIns.Range("F2").Value = "normal"
$wsIns.Range("G2").Formula = "'2012-04-02"
$wsIns.Range("H2").Value = "丁守中"
$wsIns.Range("I2").Value = 515
$wsIns.Range("J2").Value = "tmpf49e1"
$wsIns.Range("K2").Value = 105

# Row 3 (index 106)
$wsIns.Range("B3").Value = "南山人壽"
$wsIns.Range("C3").Value = "南山康福20年期終身壽險"
$wsIns.Range("D3").Value = "溫子苓"
$wsIns.Range("E3").Value = "insurance"
$wsIns.Range("F3").Value = "normal"
$wsIns.Range("G3").Formula = "'2012-04-02"
$wsIns.Range("H3").Value = "丁守中"
$wsIns.Range("I3").Value = 515
$wsIns.Range("J3").Value = "tmpf49e1"
$wsIns.Range("K3").Value = 106

# Row 4 (index 107)
$wsIns.Range("B4").Value = "南山人壽"
$wsIns.Range("C4").Value = "終身壽險"
$wsIns.Range("D4").Value = "丁守中"
$wsIns.Range("E4").Value = "insurance"
$wsIns.Range("F4").Value = "normal"
$wsIns.Range("G4").Formula = "'2012-04-02"
$wsIns.Range("H4").Value = "丁守中"
$wsIns.Range("I4").Value = 515
$wsIns.Range("J4").Value = "tmpf49e1"
$wsIns.Range("K4").Value = 107

# Row 5 (index 108)
$wsIns.Range("B5").Value = "南山人壽"
$wsIns.Range("C5").Value = "終身壽險"
$wsIns.Range("D5").Value = "溫子苓"
$wsIns.Range("E5").Value = "insurance"
$wsIns.Range("F5").Value = "normal"
$wsIns.Range("G5").Formula = "'2012-04-02"
$wsIns.Range("H5").Value = "丁守中"
$wsIns.Range("I5").Value = 515
$wsIns.Range("J5").Value = "tmpf49e1"
$wsIns.Range("K5").Value = 108

# ------------------------------------------------------------------
# 2. "債務" (debt) sheet
# ------------------------------------------------------------------
$wsDebt = $wb.Worksheets.Item("債務")

# Header row (fixes labels that previously held stray data values)
$wsDebt.Range("B1").Value = "species"
$wsDebt.Range("C1").Value = "debtor"
$wsDebt.Range("D1").Value = "owner"
$wsDebt.Range("E1").Value = "total"
$wsDebt.Range("F1").Value = "register_date"
$wsDebt.Range("G1").Value = "register_reason"
$wsDebt.Range("H1").Value = "property_category"
$wsDebt.Range("I1").Value = "category"
$wsDebt.Range("J1").Value = "date"
$wsDebt.Range("K1").Value = "legislator_name"
$wsDebt.Range("L1").Value = "legislator_id"
$wsDebt.Range("M1").Value = "source_file"
$wsDebt.Range("N1").Value = "index"

# Row 2 (index 123) - existing B2:G2 values are unchanged, append metadata
$wsDebt.Range("H2").Value = "debt"
$wsDebt.Range("I2").Value = "normal"
$wsDebt.Range("J2").Formula = "'2012-04-02"
$wsDebt.Range("K2").Value = "丁守中"
$wsDebt.Range("L2").Value = 515
$wsDebt.Range("M2").Value = "tmpf49e1"
$wsDebt.Range("N2").Value = 123
